$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "pranavk"
$ws.Range("B6").Value = "Pranav10"
